$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$nl = [char]10

# Helper to replace literal "<br/>" markers with real line breaks
function Update-CellBreaks($cell) {
    $range = $ws.Range($cell)
    $text = $range.Value2
    if ($text -ne $null) {
        $newText = $text.Replace("<br/>", $nl)
        $range.Value = $newText
    }
}

$cells = @("B20","C20","D20","B21","C21","D21","B22","C22","D22","B23","C23","D23")
foreach ($c in $cells) {
    Update-CellBreaks $c
}
